$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the three changed data values in row 4
$ws.Range("E4").Value = 6
$ws.Range("G4").Value = 3
$ws.Range("H4").Value = 13

# Move the active selection from I7 to E4
$ws.Range("E4").Select()
